$wb = $excel.ActiveWorkbook

# --- Report regenerated for handoff: update status text + timestamps ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-18 06:57:29"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-18 06:57:24"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-18 06:57:29"

# --- Widen the status columns so the longer "Ready for handoff" text fits ---
$overview.Columns.Item(5).ColumnWidth = 16.3333333333333
$overview.Columns.Item(6).ColumnWidth = 16.3333333333333
$zhcn.Columns.Item(3).ColumnWidth = 16.3333333333333
$dede.Columns.Item(3).ColumnWidth = 16.3333333333333
